{"js": "// Rename the test case id in the log table from \"316\" to \"Test_316_07_01\"\n// by adding a \"Test_\" prefix and a \"_07_01\" suffix around the existing id,\n// keeping the existing bold formatting intact.\nconst body = context.document.body;\n\nconst results = body.search(\"316\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for '316', found \" + results.items.length);\n}\n\nconst target = results.items[0];\n\n// Insert the suffix first, then the prefix, so the offsets of \"target\"\n// (the original \"316\" range) stay valid for both operations.\ntarget.insertText(\"_07_01\", Word.InsertLocation.after);\ntarget.insertText(\"Test_\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Rename the test case id in the log table from \"316\" to \"Test_316_07_01\"\n# by adding a \"Test_\" prefix and a \"_07_01\" suffix around the existing id,\n# keeping the existing bold formatting intact.\n$d = $word.ActiveDocument\n\n# Count occurrences of \"316\" first, to make sure we only touch the\n# intended Test Case ID value.\n$countRng = $d.Content\n$countFind = $countRng.Find\n$countFind.Text = \"316\"\n$countFind.MatchCase = $true\n$matchCount = 0\nwhile ($countFind.Execute()) {\n  $matchCount += 1\n  $countRng.Collapse(0)\n}\n\nif ($matchCount -ne 1) {\n  throw \"Expected exactly one match for '316', found $matchCount\"\n}\n\n# Locate the occurrence again to get its range.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"316\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n  throw \"Could not find '316' in the document.\"\n}\n\n$startPos = $rng.Start\n$endPos = $rng.End\n\n# Insert the suffix first (it sits after $startPos, so inserting it does\n# not shift $startPos), then insert the prefix at the original start.\n$d.Range($endPos, $endPos).InsertBefore(\"_07_01\")\n$d.Range($startPos, $startPos).InsertBefore(\"Test_\")\n"}
